# Auto-generated edit script applying the diff changes to before.xlsx
# Sheet order: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local life), 4=全部类型 (All types)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range('F2').Value = 233
$ws.Range('F3').Value = 423
$ws.Range('F4').Value = 163
$ws.Range('G5').Value = '不可售'
$ws.Range('F6').Value = 3890
$ws.Range('F8').Value = 2556
$ws.Range('F10').Value = 3164
$ws.Range('F11').Value = 534
$ws.Range('F14').Value = 117
$ws.Range('F15').Value = 324
$ws.Range('F16').Value = 459
$ws.Range('F17').Value = 12
$ws.Range('F18').Value = 25
$ws.Range('F19').Value = 214
$ws.Range('F22').Value = 410
$ws.Range('F23').Value = 665
$ws.Range('F24').Value = 1415
$ws.Range('F25').Value = 45
$ws.Range('F26').Value = 10
$ws.Range('F27').Value = 1306
$ws.Range('F28').Value = 134
$ws.Range('F29').Value = 154
$ws.Range('F30').Value = 29
$ws.Range('F31').Value = 5
$ws.Range('F32').Value = 62
$ws.Range('F33').Value = 4340
$ws.Range('F34').Value = 4130
$ws.Range('F35').Value = 80
$ws.Range('F36').Value = 116
$ws.Range('F38').Value = 1136
$ws.Range('F39').Value = 5
$ws.Range('F40').Value = 478
$ws.Range('F42').Value = 1318
$ws.Range('F44').Value = 132
$ws.Range('F45').Value = 106
$ws.Range('F47').Value = 64
$ws.Range('F48').Value = 63

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range('F6').Value = 2
$ws.Range('F14').Value = 9

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range('F2').Value = 1032
$ws.Range('F4').Value = 2309

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range('F2').Value = 233
$ws.Range('F3').Value = 1032
$ws.Range('F5').Value = 423
$ws.Range('F8').Value = 163
$ws.Range('C9').Value = '北京·《山丘》音乐教父 经典情歌金曲翻唱演唱会'
$ws.Range('D9').Value = '大江胡同121号2幢负1层 北京门空间 TheDoorLiveHouse'
$ws.Range('E9').Value = '2024.08.23 19:30-08.23 21:00'
$ws.Range('F9').Value = 1
$ws.Range('G9').Value = 98
$ws.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=89358'
$ws.Range('I9').Value = '//i0.hdslb.com/bfs/openplatform/202407/noqwx8Qu1721116074567.jpeg'
$ws.Range('C10').Value = '北京·喘气动漫嘉年华·暑期狂欢'
$ws.Range('D10').Value = '新风街3号 紫园·新风里'
$ws.Range('E10').Value = '2024.08.23 10:00-08.25 20:00'
$ws.Range('F10').Value = 3890
$ws.Range('G10').Value = 49
$ws.Range('H10').Value = 'https://show.bilibili.com/platform/detail.html?id=90022'
$ws.Range('I10').Value = '//i2.hdslb.com/bfs/openplatform/202408/Rb5sRto71722841653388.jpeg'
$ws.Range('B11').Value = '2024-08-24'
$ws.Range('C11').Value = ' 北京·万游引力嘉年华 配音演员赵成晨&尘霜满眸 广播剧《奕曲同工》专场见面&签售会'
$ws.Range('D11').Value = '金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心'
$ws.Range('E11').Value = '2024.08.24 11:00-08.24 17:00'
$ws.Range('F11').Value = 223
$ws.Range('G11').Value = 288
$ws.Range('H11').Value = 'https://show.bilibili.com/platform/detail.html?id=89054'
$ws.Range('I11').Value = '//i2.hdslb.com/bfs/openplatform/202407/FadWpN3x1720599868028.jpeg'
$ws.Range('C12').Value = '北京·ACY动漫游戏展1st'
$ws.Range('D12').Value = '崇文门外大街18号 北京国瑞购物中心'
$ws.Range('E12').Value = '2024.08.24 10:00-08.25 17:00'
$ws.Range('F12').Value = 2556
$ws.Range('G12').Value = 70
$ws.Range('H12').Value = 'https://show.bilibili.com/platform/detail.html?id=87851'
$ws.Range('I12').Value = '//i1.hdslb.com/bfs/openplatform/202406/jKUUlXAR1718967902236.jpeg'
$ws.Range('C13').Value = '北京·“不健全关系”专题聚会【免票活动】'
$ws.Range('D13').Value = '王府井大街88号 北京王府井银泰in88购物中心'
$ws.Range('E13').Value = '2024.08.24 14:00-08.24 18:00'
$ws.Range('F13').Value = 79
$ws.Range('G13').Value = 50
$ws.Range('H13').Value = 'https://show.bilibili.com/platform/detail.html?id=90562'
$ws.Range('I13').Value = '//i2.hdslb.com/bfs/openplatform/202408/qBI8D5Ji1723624479890.jpeg'
$ws.Range('C14').Value = '北京·万游引力夏日动漫游戏狂欢节'
$ws.Range('D14').Value = '金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心'
$ws.Range('E14').Value = '2024.08.24 10:00-08.25 17:00'
$ws.Range('F14').Value = 3164
$ws.Range('G14').Value = 75
$ws.Range('H14').Value = 'https://show.bilibili.com/platform/detail.html?id=83880'
$ws.Range('I14').Value = '//i2.hdslb.com/bfs/openplatform/202407/3EF1Am6T1720430616435.jpeg'
$ws.Range('C15').Value = '北京·最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会'
$ws.Range('D15').Value = '亮马桥路40号(近好运街) 北京世纪剧院'
$ws.Range('E15').Value = '2024.08.24 19:30-08.24 21:00'
$ws.Range('F15').Value = 19
$ws.Range('G15').Value = 238
$ws.Range('H15').Value = 'https://show.bilibili.com/platform/detail.html?id=86217'
$ws.Range('I15').Value = '//i2.hdslb.com/bfs/openplatform/202405/BDyblKrJ1716427731729.jpeg'
$ws.Range('F16').Value = 534
$ws.Range('F19').Value = 117
$ws.Range('F20').Value = 324
$ws.Range('F21').Value = 12
$ws.Range('F22').Value = 25
$ws.Range('F25').Value = 410
$ws.Range('F26').Value = 665
$ws.Range('F27').Value = 1415
$ws.Range('F28').Value = 45
$ws.Range('F29').Value = 1306
$ws.Range('F30').Value = 154
$ws.Range('F32').Value = 62
$ws.Range('F34').Value = 4340
$ws.Range('F35').Value = 4130
$ws.Range('F36').Value = 80
$ws.Range('F38').Value = 1136
$ws.Range('F39').Value = 5
$ws.Range('F43').Value = 9
$ws.Range('F45').Value = 1318
$ws.Range('F47').Value = 106

